# Add two new country sheets ("Norway" and "Poland") to the workbook,
# modeled after the existing "Croatia" sheet (same layout / column widths /
# styles as the other "Miscellaneous MZX Panels" country sheets), placed
# right after "Hungary" (the last sheet), and populate their market name /
# Jira-story cells.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Croatia")
$hungary  = $wb.Worksheets.Item("Hungary")

# --- Norway --------------------------------------------------------------
$template.Copy($null, $hungary)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B4").Value = "NGC-2931/T3085/T3078/T3084"
$norway.Range("B2").Value = "Norway Market"

# --- Poland ---------------------------------------------------------------
$template.Copy($null, $norway)
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"
$poland.Range("B4").Value = "NGC-2920/T3121/T3119/T3113"
$poland.Range("B2").Value = "Poland Market"

# Norway becomes the active/selected sheet (matches activeTab 15 -> 16).
$norway.Activate()
